# "Add files via upload" — refresh the exposure-site table with the latest
# published rows (Camberwell / Clayton South / Melbourne entries), replacing
# the previous Albert Park / Melbourne / Moorabin rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Camberwell",    "Tao Dumplings  1 Evans Place, Camberwell VIC 3124", "29/12/20 12:30pm-1:30pm",    "Case ate at restaurant",                                  "old"),
    @("Clayton South", "Metro Train - Pakenham line",                      "31/12/20 9:00pm-9:30pm",     "Case caught train from Westall Station to Flinders St.", "new"),
    @("Melbourne",     "Metro Train - Pakenham line",                      "01/01/21 3:00am-4:00am",     "Case caught train from Flinders St to Westall Station",  "old"),
    @("Melbourne",     "Metro Train - Pakenham line",                      "01/01/21 4:30am-5:00am",     "Case caught train from Flinders St to Westall Station",  "new"),
    @("Melbourne",     "Nandos  27 Elizabeth Street, Melbourne",           "01/01/2021 1:00am - 2:00am", "Case dined at venue",                                     "old"),
    @("Melbourne",     "Nandos  27 Elizabeth Street, Melbourne",           "01/01/2021 2:00am - 2:30am", "Case dined at venue",                                     "new")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $ws.Cells.Item($row, 5).Value = $entry[4]
    $row++
}

# Resize the columns to fit the new content, same as the original author did.
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(2).AutoFit() | Out-Null
$ws.Columns.Item(3).AutoFit() | Out-Null
$ws.Columns.Item(4).AutoFit() | Out-Null
$ws.Columns.Item(5).AutoFit() | Out-Null

$ws.Range("A2").Select() | Out-Null
